# Update the "Förändrad" (Changed) date column (C) for all data rows.
# The diff shows the serial date value 46061 -> 46062 (i.e. +1 day)
# for every populated row (C2:C9) on the only worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 3).End(-4162).Row  # xlUp = -4162

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    $current = $cell.Value2
    if ($current -ne $null) {
        $cell.Value2 = $current + 1
    }
}
